# Weekly data refresh ("Fruta / hortaliza, semanal"):
# two new weekly price records for Perejil (Vega Central Mapocho de Santiago)
# are inserted at rows 245-246, shifting the existing rows 245-326 down to 247-328.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the former row 245 (everything below shifts down by 2)
$ws.Rows.Item(245).Insert()
$ws.Rows.Item(245).Insert()

# Populate the first new row (245)
$ws.Cells.Item(245, 1).Value = 9
$ws.Cells.Item(245, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(245, 3).Value = "Metropolitana"
$ws.Cells.Item(245, 4).Value = 44627
$ws.Cells.Item(245, 5).Value = 13
$ws.Cells.Item(245, 6).Value = 100112044
$ws.Cells.Item(245, 7).Value = "Perejil"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 34
$ws.Cells.Item(245, 11).Value = 8000
$ws.Cells.Item(245, 12).Value = 8000
$ws.Cells.Item(245, 13).Value = 8000
$ws.Cells.Item(245, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(245, 15).Value = "Región Metropolitana"
$ws.Cells.Item(245, 16).Value = 222
$ws.Cells.Item(245, 17).Value = 36
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# Populate the second new row (246)
$ws.Cells.Item(246, 1).Value = 9
$ws.Cells.Item(246, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(246, 3).Value = "Metropolitana"
$ws.Cells.Item(246, 4).Value = 44627
$ws.Cells.Item(246, 5).Value = 13
$ws.Cells.Item(246, 6).Value = 100112044
$ws.Cells.Item(246, 7).Value = "Perejil"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 61
$ws.Cells.Item(246, 11).Value = 16000
$ws.Cells.Item(246, 12).Value = 18000
$ws.Cells.Item(246, 13).Value = 16984
$ws.Cells.Item(246, 14).Value = "`$/docena de atados"
$ws.Cells.Item(246, 15).Value = "Región Metropolitana"
$ws.Cells.Item(246, 16).Value = 5661
$ws.Cells.Item(246, 17).Value = 3
$ws.Cells.Item(246, 18).Value = "Hortaliza"

Write-Output "Inserted 2 rows; new dimension should be A1:R328"
